$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data in row 8 (ticket sales row) - reset to 0
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

# Resize columns B and C to fit the new, wider content
# (ColumnWidth is in characters; Excel stores the saved XML width in whole
# pixels = round(chars*6)+5, all /6 -- so we dial in the character width that
# lands on the desired stored width of 23.5 / 37.625 after that padding.)
$ws.Columns.Item(2).ColumnWidth = 22.67
$ws.Columns.Item(3).ColumnWidth = 36.83

# Move the active selection to H14
$ws.Range("H14").Select()
